$d = $word.ActiveDocument

# 1. Update "Curso (semestre ideal)" line to add EQD (3)
$d.Content.Find.Execute("Curso (semestre ideal): EQN (3)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Curso (semestre ideal): EQD (3), EQN (3)", 2)

# 2. Remove the "Requisitos" heading paragraph and the following
#    "LOQ4073 -  Química Geral II  (Requisito fraco)" bullet paragraph.
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    if ($paras.Item($i).Range.Text.TrimEnd([char]13, [char]7) -eq "Requisitos") {
        $startPara = $paras.Item($i)
        $endPara = $paras.Item($i + 1)
        $range = $d.Range($startPara.Range.Start, $endPara.Range.End)
        $range.Delete()
        break
    }
}
